$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell carrying the default (unstyled, General-format) style
# used to reset style after forcing text entry via a leading apostrophe.
$defaultStyle = $ws.Range("D4").Style

function Set-TextValue($addr, $value) {
    # Force the literal to be stored as TEXT (matches the workbook's
    # existing inline-string cells) by prefixing with an apostrophe,
    # then restore the plain default style so no stray "quote prefix"
    # cell format lingers on the cell.
    $ws.Range($addr).Formula = "`'" + $value
    $ws.Range($addr).Style = $defaultStyle
}

# Row 2
Set-TextValue "D2" "73.101.51"
$ws.Range("E2").Value = "  +1.47%  "

# Row 3
Set-TextValue "D3" "4.048.36"
$ws.Range("E3").Value = "  +0.84%  "

# Row 4
$ws.Range("E4").Value = "  -0.20%  "

# Row 5
Set-TextValue "D5" "594.42"
$ws.Range("E5").Value = "  +12.36%  "

# Row 6
Set-TextValue "D6" "153.98"
$ws.Range("E6").Value = "  +2.59%  "

# Row 7
Set-TextValue "D7" "0.690"
$ws.Range("E7").Value = "  -1.06%  "

# Row 8
Set-TextValue "D8" "0.999"
$ws.Range("E8").Value = "  -0.10%  "

# Row 9
Set-TextValue "D9" "0.764"
$ws.Range("E9").Value = "  +2.33%  "

# Row 10
Set-TextValue "D10" "0.171"
$ws.Range("E10").Value = "  +0.39%  "

# Row 11
Set-TextValue "D11" "53.88"

# Row 12
$ws.Range("E12").Value = "  -0.27%  "

# Row 13
Set-TextValue "D13" "11.09"
$ws.Range("E13").Value = "  +4.83%  "

# Row 14
Set-TextValue "D14" "4.701.30"
$ws.Range("E14").Value = "  +0.69%  "

# Row 15
Set-TextValue "D15" "4.049.13"
$ws.Range("E15").Value = "  +0.68%  "

# Row 16
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D16" "1.26"
$ws.Range("E16").Value = "  +6.17%  "

# Row 17
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D17" "14.31"
$ws.Range("E17").Value = "  +1.95%  "

# Row 18
Set-TextValue "D18" "20.72"
$ws.Range("E18").Value = "  +1.35%  "

# Row 19
$ws.Range("E19").Value = "  -0.53%  "

# Row 20
Set-TextValue "D20" "73.014.34"
$ws.Range("E20").Value = "  +1.31%  "

# Row 21
Set-TextValue "D21" "444.73"
$ws.Range("E21").Value = "  +4.21%  "

# Row 22
Set-TextValue "D22" "4.76"
$ws.Range("E22").Value = "  +13.84%  "

# Row 23
Set-TextValue "D23" "97.59"
$ws.Range("E23").Value = "  +0.55%  "

# Row 24
$ws.Range("E24").Value = "  +1.96%  "

# Row 25
Set-TextValue "D25" "14.38"
$ws.Range("E25").Value = "  +1.94%  "

# Row 26
Set-TextValue "D26" "4.38"
$ws.Range("E26").Value = "  +22.57%  "

# Row 27
Set-TextValue "D27" "11.35"
$ws.Range("E27").Value = "  +1.76%  "

# Row 28
$ws.Range("E28").Value = "  +2.04%  "

# Row 29
Set-TextValue "D29" "5.93"
$ws.Range("E29").Value = "  +1.59%  "

# Row 30
Set-TextValue "D30" "36.97"
$ws.Range("E30").Value = "  +1.10%  "

# Row 31
Set-TextValue "D31" "8.02"
$ws.Range("E31").Value = "  +14.02%  "

# Row 32
$ws.Range("E32").Value = "  +4.80%  "

# Row 33
$ws.Range("E33").Value = "  +2.97%  "

# Row 34
Set-TextValue "D34" "687.73"
$ws.Range("E34").Value = "  +1.67%  "

# Row 35
Set-TextValue "D35" "49.22"
$ws.Range("E35").Value = "  +11.86%  "

# Row 36
Set-TextValue "D36" "71.14"
$ws.Range("E36").Value = "  +7.91%  "

# Row 37
Set-TextValue "D37" "0.450"
$ws.Range("E37").Value = "  +2.13%  "

# Row 38
Set-TextValue "D38" "0.0₃0880"
$ws.Range("E38").Value = "  +6.22%  "

# Row 39
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D39" "0.149"
$ws.Range("E39").Value = "  -1.57%  "

# Row 40
Set-TextValue "D40" "11.31"
$ws.Range("E40").Value = "  +16.57%  "

# Row 41
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextValue "D41" "3.38"
$ws.Range("E41").Value = "  +0.08%  "

# Row 42
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D42" "3.37"
$ws.Range("E42").Value = "  +5.47%  "

# Row 43
$ws.Range("E43").Value = "  +0.07%  "

# Row 44
$ws.Range("E44").Value = "  +2.47%  "

# Row 45
$ws.Range("E45").Value = "  +0.23%  "

# Row 46
$ws.Range("E46").Value = "  +1.57%  "

# Row 47
$ws.Range("E47").Value = "  +3.06%  "

# Row 48
$ws.Range("B48").Value = "LidoDAOToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D48" "3.55"
$ws.Range("E48").Value = "  +9.03%  "

# Row 49
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D49" "3.37"
$ws.Range("E49").Value = "  -1.55%  "

# Row 50
$ws.Range("E50").Value = "  +2.19%  "

# Row 51
$ws.Range("E51").Value = "  +9.54%  "

Write-Output "Applied cryptos list update"
